$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in D4 and D5
$ws.Range("D4").Value = 0.01749518539312528
$ws.Range("D5").Value = 0.006017676599826928

# Add new row 6 with DWA data, copying the formatting from A5 (name-column style)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").Value = "DWA"

$ws.Range("B6").Value = 65
$ws.Range("C6").Value = 0.6445621764871753
$ws.Range("D6").Value = 0.0149623279218316
$ws.Range("E6").Value = 0.6308837203990878
